# Re-run of the Leidingweerstand model-coefficient fit (new `use_lims`
# argument + Aquathermie plots). Refreshes the "gewijzigd" timestamp on
# every coefficient sheet, and on Q100 a new data point was added to the
# regression (row shifts down, coefficients recomputed).

$wb = $excel.ActiveWorkbook

# ---- Q100: new data point inserted, coefficients + timestamps refreshed ----
$ws = $wb.Worksheets.Item("Q100")
$ws.Range("C2").Value = [double]"-9.580804480917463e-06"
$ws.Range("E2").Value = 45096.5062931713
$ws.Range("C3").Value = [double]"-1.196993856286746e-05"
$ws.Range("E3").Value = 45096.50629331019
$ws.Range("B4").Value = 43467
$ws.Range("C4").Value = [double]"-1.235191816021002e-05"
$ws.Range("E4").Value = 45096.5062931713
$ws.Range("A5").Value = 3
$ws.Range("A5").Style = $ws.Range("A4").Style
$ws.Range("B5").Value = 44443.5
$ws.Range("B5").Style = $ws.Range("B4").Style
$ws.Range("C5").Value = [double]"-1.21000272712855e-05"
$ws.Range("D5").Value = [double]"-1.2032909e-09"
$ws.Range("E5").Value = 45096.5062931713
$ws.Range("E5").Style = $ws.Range("E4").Style

# ---- IK106: only the "gewijzigd" timestamps move ----
$ws = $wb.Worksheets.Item("IK106")
$ws.Range("E2").Value = 45096.51158341607
$ws.Range("E3").Value = 45096.51158341607
$ws.Range("E4").Value = 45096.51158357813
$ws.Range("E5").Value = 45096.51158341607

# ---- Q200 ----
$ws = $wb.Worksheets.Item("Q200")
$ws.Range("E2").Value = 45096.50654289352
$ws.Range("E3").Value = 45096.50654289352
$ws.Range("E4").Value = 45096.50654289352

# ---- Q300 ----
$ws = $wb.Worksheets.Item("Q300")
$ws.Range("E2").Value = 45096.50686523148
$ws.Range("E3").Value = 45096.50686537037
$ws.Range("E4").Value = 45096.50686537037
$ws.Range("E5").Value = 45096.50686523148

# ---- Q400 ----
$ws = $wb.Worksheets.Item("Q400")
$ws.Range("E2").Value = 45096.50712299768
$ws.Range("E3").Value = 45096.50712299768
$ws.Range("E4").Value = 45096.50712299768
$ws.Range("E5").Value = 45096.50712299768

# ---- Q500 ----
$ws = $wb.Worksheets.Item("Q500")
$ws.Range("E2").Value = 45096.50749451389
$ws.Range("E3").Value = 45096.50749466435
$ws.Range("E4").Value = 45096.50749466435
$ws.Range("E5").Value = 45096.50749451389
$ws.Range("E6").Value = 45096.50749451389

# ---- Q600 ----
$ws = $wb.Worksheets.Item("Q600")
$ws.Range("E2").Value = 45096.50780311343
$ws.Range("E3").Value = 45096.50780325232
$ws.Range("E4").Value = 45096.50780311343
$ws.Range("E5").Value = 45096.50780311343

# ---- P100 ----
$ws = $wb.Worksheets.Item("P100")
$ws.Range("E2").Value = 45096.5080496412
$ws.Range("E3").Value = 45096.5080496412
$ws.Range("E4").Value = 45096.5080496412
$ws.Range("E5").Value = 45096.5080496412

# ---- P200 ----
$ws = $wb.Worksheets.Item("P200")
$ws.Range("E2").Value = 45096.50840765046
$ws.Range("E3").Value = 45096.50840765046
$ws.Range("E4").Value = 45096.50840782408
$ws.Range("E5").Value = 45096.50840765046
$ws.Range("E6").Value = 45096.50840765046

# ---- P300 ----
$ws = $wb.Worksheets.Item("P300")
$ws.Range("E2").Value = 45096.50893650463
$ws.Range("E3").Value = 45096.50893666667
$ws.Range("E4").Value = 45096.50893666667
$ws.Range("E5").Value = 45096.50893666667
$ws.Range("E6").Value = 45096.50893650463
$ws.Range("E7").Value = 45096.50893650463
$ws.Range("E8").Value = 45096.50893650463

# ---- P400 ----
$ws = $wb.Worksheets.Item("P400")
$ws.Range("E2").Value = 45096.50927262731
$ws.Range("E3").Value = 45096.50927262731
$ws.Range("E4").Value = 45096.5092727662
$ws.Range("E5").Value = 45096.50927262731
$ws.Range("E6").Value = 45096.50927262731

# ---- P500 ----
$ws = $wb.Worksheets.Item("P500")
$ws.Range("E2").Value = 45096.50959040509
$ws.Range("E3").Value = 45096.50959050926
$ws.Range("E4").Value = 45096.50959050926
$ws.Range("E5").Value = 45096.50959040509

# ---- P600 ----
$ws = $wb.Worksheets.Item("P600")
$ws.Range("E2").Value = 45096.509836875
$ws.Range("E3").Value = 45096.509836875
$ws.Range("E4").Value = 45096.509836875
$ws.Range("E5").Value = 45096.509836875

# ---- IK91 ----
$ws = $wb.Worksheets.Item("IK91")
$ws.Range("E2").Value = 45096.51002239584
$ws.Range("E3").Value = 45096.51002239584
$ws.Range("E4").Value = 45096.51002239584

# ---- IK92 ----
$ws = $wb.Worksheets.Item("IK92")
$ws.Range("E2").Value = 45096.51011625
$ws.Range("E3").Value = 45096.51011625
$ws.Range("E4").Value = 45096.51011625

# ---- IK93 ----
$ws = $wb.Worksheets.Item("IK93")
$ws.Range("E2").Value = 45096.51023998843
$ws.Range("E3").Value = 45096.51023998843
$ws.Range("E4").Value = 45096.51023998843
$ws.Range("E5").Value = 45096.51023998843
$ws.Range("E6").Value = 45096.51023998843

# ---- IK94 ----
$ws = $wb.Worksheets.Item("IK94")
$ws.Range("E2").Value = 45096.51034395833
$ws.Range("E3").Value = 45096.51034395833
$ws.Range("E4").Value = 45096.51034395833
$ws.Range("E5").Value = 45096.51034395833

# ---- IK95 ----
$ws = $wb.Worksheets.Item("IK95")
$ws.Range("E2").Value = 45096.510441875
$ws.Range("E3").Value = 45096.510441875
$ws.Range("E4").Value = 45096.510441875
$ws.Range("E5").Value = 45096.510441875

# ---- IK96 ----
$ws = $wb.Worksheets.Item("IK96")
$ws.Range("E2").Value = 45096.51056390046
$ws.Range("E3").Value = 45096.51056390046
$ws.Range("E4").Value = 45096.51056390046

# ---- IK101 ----
$ws = $wb.Worksheets.Item("IK101")
$ws.Range("E2").Value = 45096.51066295139
$ws.Range("E3").Value = 45096.51066295139
$ws.Range("E4").Value = 45096.51066295139

# ---- IK102 ----
$ws = $wb.Worksheets.Item("IK102")
$ws.Range("E2").Value = 45096.51076696759
$ws.Range("E3").Value = 45096.51076696759
$ws.Range("E4").Value = 45096.51076696759
$ws.Range("E5").Value = 45096.51076696759

# ---- IK103 ----
$ws = $wb.Worksheets.Item("IK103")
$ws.Range("E2").Value = 45096.51097048611
$ws.Range("E3").Value = 45096.51097072916
$ws.Range("E4").Value = 45096.51097048611
$ws.Range("E5").Value = 45096.51097048611
$ws.Range("E6").Value = 45096.51097072916

# ---- IK104 ----
$ws = $wb.Worksheets.Item("IK104")
$ws.Range("E2").Value = 45096.51120732639
$ws.Range("E3").Value = 45096.5112075463
$ws.Range("E4").Value = 45096.51120732639
$ws.Range("E5").Value = 45096.51120732639

# ---- IK105 ----
$ws = $wb.Worksheets.Item("IK105")
$ws.Range("E2").Value = 45096.51134891203
$ws.Range("E3").Value = 45096.51134891203
$ws.Range("E4").Value = 45096.51134891203
$ws.Range("E5").Value = 45096.51134891203
